$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp shown in row 1
$ws.Range("A1").Value = "Datos actualizados a 19 de Marzo de 2020 a las 21:14"

# Insert a new row for "Cataluna*" right after Madrid (row 4), carrying the
# same totals as Cataluña (which shifts from row 5 down to row 6).
$ws.Rows(5).Insert()
$ws.Range("A5").Value = "Cataluna*"
$ws.Range("B5").Value = 3270
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 3185
$ws.Range("E5").Value = 82

# "Illes Balears" and "Illes Balears*" swap order (rows 25/26 after the insert above).
$ws.Range("A25").Value = "Illes Balears"
$ws.Range("A26").Value = "Illes Balears*"

# "Ceuta" and "La Palma" swap order (rows 58/59 after the insert above).
$ws.Range("A58").Value = "Ceuta"
$ws.Range("A59").Value = "La Palma"
